# Applies the "missing_axis_name" + "narrow_2d" sheets addition described
# by the diff: two new worksheets are inserted right before "2d_classic"
# (pushing it, "unsorted" and "int_labels" one slot to the right), and the
# workbook's active tab moves to the newly inserted "narrow_2d" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the two new sheets in the right spot / right order so that
#    sheetId allocation + r:id allocation matches the target (sheetId 9
#    for missing_axis_name, sheetId 10 for narrow_2d; rId5/rId6).
#    NOTE: Worksheets.Item(...) is a *live*, position-bound reference in
#    this host -- it must be re-fetched by name after every insertion,
#    otherwise it silently re-resolves to whatever sheet now sits at the
#    old index.
# ---------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item("2d_classic")
$wsMissingAxis = $wb.Worksheets.Add($refSheet)
$wsMissingAxis.Name = "missing_axis_name"

$refSheet2 = $wb.Worksheets.Item("2d_classic")
$wsNarrow2d = $wb.Worksheets.Add($refSheet2)
$wsNarrow2d.Name = "narrow_2d"

# ---------------------------------------------------------------------
# 2) Populate "missing_axis_name" (A1:D5)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("missing_axis_name")

$ws.Range("A1").Value = "a"
$ws.Range("B1").Value = "b"
$ws.Range("C1").Value = "c0"
$ws.Range("D1").Value = "c1"

$ws.Range("A2").Value = "a0"
$ws.Range("B2").Value = "b0"
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 1

$ws.Range("A3").Value = "a0"
$ws.Range("B3").Value = "b1"
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 3

$ws.Range("A4").Value = "a1"
$ws.Range("B4").Value = "b0"
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 5

$ws.Range("A5").Value = "a1"
$ws.Range("B5").Value = "b1"
$ws.Range("C5").Value = 6
$ws.Range("D5").Value = 7

# column A (the repeated row-axis labels) is vertically centered and
# uses an explicit (non-themed) Calibri font -- mirrors the style added
# to xl/styles.xml (cellXfs[1]).
$ws.Range("A1:A5").Font.Name = "Calibri"
$ws.Range("A1:A5").VerticalAlignment = -4108

[void]$ws.Range("G18").Select()

# ---------------------------------------------------------------------
# 3) Populate "narrow_2d" (A1:C7)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("narrow_2d")

$ws2.Range("A1").Value = "a"
$ws2.Range("B1").Value = "b"
$ws2.Range("C1").Value = "value"

$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = "b0"
$ws2.Range("C2").Value = 0

$ws2.Range("A3").Value = 1
$ws2.Range("B3").Value = "b1"
$ws2.Range("C3").Value = 1

$ws2.Range("A4").Value = 2
$ws2.Range("B4").Value = "b0"
$ws2.Range("C4").Value = 2

$ws2.Range("A5").Value = 2
$ws2.Range("B5").Value = "b1"
$ws2.Range("C5").Value = 3

$ws2.Range("A6").Value = 3
$ws2.Range("B6").Value = "b0"
$ws2.Range("C6").Value = 4

$ws2.Range("A7").Value = 3
$ws2.Range("B7").Value = "b1"
$ws2.Range("C7").Value = 5

[void]$ws2.Range("F17").Select()

# "narrow_2d" ends up the active tab (activeTab goes 6 -> 5, i.e. the
# sheet now sitting at position 5 once the other sheets shifted right).
$ws2.Activate()
